# semana 5 de 2026
# Adds week-5 column (H) data to the weekly IRA extract sheet and inserts
# a newly-reporting UPGD ("CENTRO DE SALUD SAN NICOLAS") as row 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new UPGD row at position 10 (pushes old rows 10-47 -> 11-48) ---
$ws.Rows(10).Insert()

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "6600100332"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "05"
$ws.Range("C10").Value = "CENTRO DE SALUD SAN NICOLAS"
$ws.Range("H10").Value = 2

# --- New header cell for week 5, matching the style of the other header cells ---
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").Value = "5"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108

# --- Week-5 (column H) counts for every UPGD row ---
$ws.Range("H2").Value = 20
$ws.Range("H3").Value = 95
$ws.Range("H4").Value = 3
$ws.Range("H5").Value = 59
$ws.Range("H6").Value = 27
$ws.Range("H7").Value = 40
$ws.Range("H8").Value = 2
$ws.Range("H9").Value = 1
# H10 already set above (new row)
# H11 (CENTRO DE SALUD SAN CAMILO) and H12 (CENTRO DE SALUD VILLASANTANA) have no week-5 data
$ws.Range("H13").Value = 9
$ws.Range("H14").Value = 4
$ws.Range("H15").Value = 4
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 2
$ws.Range("H18").Value = 10
$ws.Range("H19").Value = 34
# H20 (CENTRO DE SALUD CASA DEL ABUELO) has no week-5 data
$ws.Range("H21").Value = 171
$ws.Range("H22").Value = 0
$ws.Range("H23").Value = 47
$ws.Range("H24").Value = 4
$ws.Range("H25").Value = 40
$ws.Range("H26").Value = 6
$ws.Range("H27").Value = 14
$ws.Range("H28").Value = 34
# H29 (IPS CENTRO DE MEDICINA INTEGRATIVA SAS) has no week-5 data
$ws.Range("H30").Value = 9
$ws.Range("H31").Value = 99
$ws.Range("H32").Value = 48
$ws.Range("H33").Value = 8
$ws.Range("H34").Value = 216
$ws.Range("H35").Value = 93
$ws.Range("H36").Value = 127
$ws.Range("H37").Value = 4
$ws.Range("H38").Value = 119
$ws.Range("H39").Value = 6
$ws.Range("H40").Value = 0
$ws.Range("H41").Value = 7
$ws.Range("H42").Value = 6
$ws.Range("H43").Value = 0
$ws.Range("H44").Value = 3
$ws.Range("H45").Value = 0
$ws.Range("H46").Value = 10
$ws.Range("H47").Value = 88
$ws.Range("H48").Value = 25
